$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that are fully removed in the target (rows 2-3 Action2/x2/Character2 cols, and rows 10-13 Vocal/BGImage cols) ---
$ws.Range("J2:L3").Clear()
$ws.Range("D10:E13").Clear()
# Row 8's Action2/x2/Character2 move to row 4; keep the (now empty) styled placeholder cells on row 8
$ws.Range("J8:L8").ClearContents()

# --- Update row 2 ---
$ws.Range("A2").Value = "Yao"
$ws.Range("B2").Value = "Sir, there’s a note on the desk."
$ws.Range("C2").Value = "Yao-Regular"
$ws.Range("D2").Value = "DialogueVocal"
$ws.Range("E2").Value = "StudyInvestigate"
$ws.Range("F2").Value = "Suspicious"

# --- Update row 3 ---
$ws.Range("A3").Value = "Dee"
$ws.Range("B3").Value = "What is this......?"
$ws.Range("C3").Value = "Dee-Determined"
$ws.Range("D3").Value = "DialogueVocal"
$ws.Range("E3").Value = "StudyInvestigate"

# --- Update row 4 ---
$ws.Range("B4").Value = " <color=#00CC00>(“Come outside the manor at the second quarter of the Xu hour.”)</color>"
$ws.Range("D4").Value = "DialogueVocal"
$ws.Range("E4").Value = "StudyInvestigate"
$ws.Range("J4").Value = "appearAt"
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = "Evi-Note"

# --- Update row 5 ---
$ws.Range("A5").Value = "Yao"
$ws.Range("B5").Value = "This must be the note the culprit sent to the Lord, asking to meet outside the manor 7.30 PM."
$ws.Range("C5").Value = "Yao-Regular"
$ws.Range("D5").Value = "DialogueVocal"
$ws.Range("E5").Value = "StudyInvestigate"

# --- Update row 6 ---
$ws.Range("A6").Value = "Dee"
$ws.Range("B6").Value = "The handwriting looks awkward——whoever wrote it must have used their left hand."
$ws.Range("C6").Value = "Dee-Thinking2"
$ws.Range("D6").Value = "DialogueVocal"
$ws.Range("E6").Value = "StudyInvestigate"

# --- Update row 7 (add D7/E7) ---
$ws.Range("B7").Value = " <color=#00CC00>(The texture and style of the paper also differ from ordinary letter paper.)</color>"
$ws.Range("D7").Value = "DialogueVocal"
$ws.Range("E7").Value = "StudyInvestigate"

# --- Update row 8 ---
$ws.Range("B8").Value = " <color=#00CC00>(Something about this note feels...... off.)</color>"
$ws.Range("D8").Value = "DialogueVocal"
$ws.Range("E8").Value = "StudyInvestigate"

# --- Update row 9 (add A9/C9) ---
$ws.Range("A9").Value = "Dee"
$ws.Range("B9").Value = "Anyway, let’s hold on to it for now."
$ws.Range("C9").Value = "Dee-Thinking2"
$ws.Range("D9").Value = "DialogueVocal"
$ws.Range("E9").Value = "StudyInvestigate"

# --- Update row 10 ---
$ws.Range("A10").Value = "Investigate"
$ws.Range("B10").Value = "Desk"
$ws.Range("C10").Value = "Desk"

# --- Update row 11 ---
$ws.Range("B11").Value = "Paper"
$ws.Range("C11").Value = "Paper"

# --- Update row 12 ---
$ws.Range("B12").Value = "Book"
$ws.Range("C12").Value = "Book"

# --- Update row 13 ---
$ws.Range("B13").Value = "End Investigation"
$ws.Range("C13").Value = "StoryScript13"

# --- Fix row heights that the auto-wrap heuristic would get wrong (shrinking cases) ---
$ws.Rows.Item(5).RowHeight = 34
$ws.Rows.Item(8).RowHeight = 34

# --- Selection ---
$ws.Range("B18").Select()
